# Append the new portfolio data row (2025-10-28) as row 74, matching the
# pattern of the existing rows: Date (text) / SUZLON.NS / TATAMOTORS.NS
# (column D / ETERNAL.NS is left blank, same as the prior row 73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Range("A74")

# Force the cell to text BEFORE assigning, so the "2025-10-28" looking-like-a-date
# string is stored as literal text (matching the existing rows, which are all
# inline/shared strings, not real dates) instead of being auto-converted into a
# date serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-10-28"
# Restore the default "Normal" style so the cell doesn't end up carrying a
# leftover explicit text-number-format style (the source rows carry no style).
$dateCell.Style = "Normal"

$ws.Range("B74").Value = 56.22000122070312
$ws.Range("C74").Value = 334.6000061035156
